# "Cambio en la forma de generar, mas optimizado"
#
# Hoja1!K3 is a single "level" input that every pricing/cost formula on the
# sheet (columns B, C, D, H, L for rows 5-9, via absolute refs to $K$3)
# depends on. The author simply lowered that input from 3 to 1 (cheaper /
# more optimized starting level) and re-saved; Excel then recalculated all
# the dependent formulas automatically, which is why every cached <v> in
# the diff for B5:L9 changed even though none of those formulas themselves
# were edited.
#
# The saved file also shows the cursor left on D17 (instead of K3) and the
# sheet scrolled a bit - incidental window/view state from that editing
# session rather than a data change, but we replicate the selection too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Core data edit: the "level" parameter everything else derives from.
$ws.Range("K3").Value = 1

# Leave the selection where the author left it when they saved.
$ws.Range("D17").Select() | Out-Null

# Force a full recalculation so every formula that references $K$3 picks up
# its new cached value (the runtime also auto-recalcs after the script, but
# doing it explicitly here keeps the result deterministic either way).
$excel.CalculateFull()
